$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '43.980.34'
Set-TextValue 'E2' '  +0.55%  '
Set-TextValue 'D3' '2.260.74'
Set-TextValue 'E3' '  -0.23%  '
Set-TextValue 'E4' '  +0.09%  '
Set-TextValue 'D5' '230.74'
Set-TextValue 'E5' '  +0.01%  '
Set-TextValue 'D6' '0.633'
Set-TextValue 'E6' '  +1.06%  '
Set-TextValue 'D7' '64.40'
Set-TextValue 'E7' '  +5.18%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.451'
Set-TextValue 'E9' '  +7.30%  '
Set-TextValue 'E10' '  +6.89%  '
Set-TextValue 'D11' '57.18'
Set-TextValue 'E11' '  -1.52%  '
Set-TextValue 'D12' '27.00'
Set-TextValue 'E12' '  +14.56%  '
Set-TextValue 'E13' '  +1.99%  '
Set-TextValue 'D14' '2.600.61'
Set-TextValue 'E14' '  -0.11%  '
Set-TextValue 'D15' '15.72'
Set-TextValue 'E15' '  +0.77%  '
Set-TextValue 'D16' '6.08'
Set-TextValue 'E16' '  +4.90%  '
Set-TextValue 'E17' '  +3.08%  '
Set-TextValue 'D18' '2.245.40'
Set-TextValue 'E18' '  -0.70%  '
Set-TextValue 'D19' '43.917.57'
Set-TextValue 'E19' '  +2.40%  '
Set-TextValue 'D20' '0.0000101'
Set-TextValue 'E20' '  +7.68%  '
Set-TextValue 'D21' '73.47'
Set-TextValue 'E21' '  +0.71%  '
Set-TextValue 'E22' '  -2.17%  '
Set-TextValue 'D23' '251.46'
Set-TextValue 'E23' '  -0.73%  '
Set-TextValue 'D24' '1.00'
Set-TextValue 'E24' '  +0.00%  '
Set-TextValue 'D25' '2.44'
Set-TextValue 'E25' '  -4.02%  '
Set-TextValue 'D26' '10.12'
Set-TextValue 'E26' '  +2.73%  '
Set-TextValue 'E27' '  -2.04%  '
Set-TextValue 'E28' '  +20.99%  '
Set-TextValue 'D29' '170.99'
Set-TextValue 'E29' '  +0.38%  '
Set-TextValue 'E30' '  -0.12%  '
Set-TextValue 'D31' '20.91'
Set-TextValue 'E31' '  +2.07%  '
Set-TextValue 'E32' '  -3.49%  '
Set-TextValue 'D33' '0.124'
Set-TextValue 'E33' '  +1.60%  '
Set-TextValue 'D34' '0.0706'
Set-TextValue 'E34' '  +6.90%  '
Set-TextValue 'E35' '  -0.05%  '
Set-TextValue 'E36' '  -3.63%  '
Set-TextValue 'D37' '3.78'
Set-TextValue 'E37' '  +5.02%  '
Set-TextValue 'D38' '6.48'
Set-TextValue 'E38' '  +0.30%  '
Set-TextValue 'D39' '2.30'
Set-TextValue 'E39' '  -3.89%  '
Set-TextValue 'E40' '  +3.71%  '
Set-TextValue 'E41' '  -0.04%  '
Set-TextValue 'D42' '0.000222'
Set-TextValue 'E42' '  -2.59%  '
Set-TextValue 'D43' '0.0974'
Set-TextValue 'E43' '  -1.87%  '
Set-TextValue 'D44' '17.32'
Set-TextValue 'E44' '  +3.95%  '
Set-TextValue 'D46' '98.06'
Set-TextValue 'E46' '  -0.18%  '
Set-TextValue 'D47' '1.20'
Set-TextValue 'E47' '  -0.56%  '
Set-TextValue 'D48' '4.42'
Set-TextValue 'E48' '  -1.93%  '
Set-TextValue 'D49' '10.26'
Set-TextValue 'E49' '  +7.40%  '
Set-TextValue 'E50' '  +5.36%  '
Set-TextValue 'D51' '1.440.30'
Set-TextValue 'E51' '  -2.24%  '
